$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column, new text value (all target cells are plain text,
# matching the original inlineStr cell type in the workbook). The leading
# comma on each row keeps it as a single 3-element array element instead of
# being unrolled into the outer $changes array.
$changes = @(
    ,@(2, 4, '35.373.95')
    ,@(2, 5, '  -0.32%  ')
    ,@(3, 4, '1.904.86')
    ,@(3, 5, '  +0.19%  ')
    ,@(4, 5, '  -0.11%  ')
    ,@(5, 2, 'XRP')
    ,@(5, 3, 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp')
    ,@(5, 4, '0.700')
    ,@(5, 5, '  +10.61%  ')
    ,@(6, 2, 'BNB')
    ,@(6, 3, 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb')
    ,@(6, 4, '246.89')
    ,@(6, 5, '  +0.73%  ')
    ,@(7, 5, '  -0.06%  ')
    ,@(8, 4, '40.82')
    ,@(8, 5, '  -2.74%  ')
    ,@(9, 4, '0.352')
    ,@(9, 5, '  +3.69%  ')
    ,@(10, 4, '52.53')
    ,@(10, 5, '  +8.05%  ')
    ,@(11, 4, '0.0726')
    ,@(11, 5, '  +3.33%  ')
    ,@(12, 5, '  -0.77%  ')
    ,@(13, 4, '2.178.95')
    ,@(13, 5, '  +0.05%  ')
    ,@(14, 4, '12.62')
    ,@(14, 5, '  +2.23%  ')
    ,@(15, 4, '0.712')
    ,@(15, 5, '  +3.32%  ')
    ,@(16, 2, 'Polkadot')
    ,@(16, 3, 'https://coinranking.com/coin/25W7FG7om+polkadot-dot')
    ,@(16, 4, '4.88')
    ,@(16, 5, '  +0.49%  ')
    ,@(17, 2, 'WrappedEther')
    ,@(17, 3, 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth')
    ,@(17, 4, '1.893.52')
    ,@(17, 5, '  -0.03%  ')
    ,@(18, 4, '35.337.19')
    ,@(18, 5, '  -0.40%  ')
    ,@(19, 4, '72.70')
    ,@(19, 5, '  +1.13%  ')
    ,@(20, 4, '0.0₃0824')
    ,@(20, 5, '  +0.39%  ')
    ,@(21, 4, '241.68')
    ,@(21, 5, '  -0.52%  ')
    ,@(22, 4, '12.91')
    ,@(22, 5, '  +3.37%  ')
    ,@(23, 4, '5.10')
    ,@(23, 5, '  +5.01%  ')
    ,@(24, 5, '  -0.06%  ')
    ,@(25, 5, '  +1.07%  ')
    ,@(26, 5, '  +6.21%  ')
    ,@(27, 4, '168.71')
    ,@(27, 5, '  -1.89%  ')
    ,@(28, 4, '8.66')
    ,@(28, 5, '  +1.13%  ')
    ,@(29, 4, '18.87')
    ,@(29, 5, '  +5.10%  ')
    ,@(30, 5, '  +5.03%  ')
    ,@(32, 4, '4.26')
    ,@(32, 5, '  +3.85%  ')
    ,@(33, 4, '0.0574')
    ,@(33, 5, '  +1.03%  ')
    ,@(34, 5, '  +0.01%  ')
    ,@(35, 5, '  +7.20%  ')
    ,@(36, 4, '4.17')
    ,@(36, 5, '  +0.00%  ')
    ,@(37, 4, '0.915')
    ,@(37, 5, '  -4.99%  ')
    ,@(38, 5, '  +7.95%  ')
    ,@(39, 4, '2.03')
    ,@(39, 5, '  +0.61%  ')
    ,@(40, 4, '98.31')
    ,@(40, 5, '  +8.47%  ')
    ,@(41, 2, 'InjectiveProtocol')
    ,@(41, 3, 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj')
    ,@(41, 4, '16.59')
    ,@(41, 5, '  +5.06%  ')
    ,@(42, 2, 'ARBITRUM')
    ,@(42, 3, 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb')
    ,@(42, 4, '1.10')
    ,@(42, 5, '  -0.13%  ')
    ,@(43, 2, 'Kaspa')
    ,@(43, 3, 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas')
    ,@(43, 4, '0.0651')
    ,@(43, 5, '  +7.53%  ')
    ,@(44, 5, '  +1.62%  ')
    ,@(45, 4, '1.360.42')
    ,@(45, 5, '  +0.77%  ')
    ,@(46, 4, '2.41')
    ,@(46, 5, '  +2.37%  ')
    ,@(47, 5, '  +0.04%  ')
    ,@(48, 5, '  +1.18%  ')
    ,@(49, 2, 'MultiversX')
    ,@(49, 3, 'https://coinranking.com/coin/omwkOTglq+multiversx-egld')
    ,@(49, 4, '45.95')
    ,@(49, 5, '  -6.23%  ')
    ,@(50, 2, 'Gas')
    ,@(50, 3, 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas')
    ,@(50, 4, '12.31')
    ,@(50, 5, '  -6.09%  ')
    ,@(51, 4, '6.53')
    ,@(51, 5, '  -1.31%  ')
)

foreach ($change in $changes) {
    $r = $change[0]
    $c = $change[1]
    $v = $change[2]
    $cell = $ws.Cells.Item($r, $c)
    # Force text format so numeric-looking strings (e.g. '1.10', '0.700')
    # are stored as text, not coerced to numbers, matching the source data.
    $cell.NumberFormat = '@'
    $cell.Value = $v
    # Reset style back to the sheet default so no stray per-cell formatting
    # is introduced beyond the text content change.
    $cell.Style = 'Normal'
}

